$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.389.63'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.790.76'
$ws.Range('E3').Value = '  -1.59%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''224.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('E6').Value = '  -3.71%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''32.71'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.95%  '
$ws.Range('D9').Value = '''0.284'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').Value = '''0.0932'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '2.047.96'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').Value = '''11.06'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +8.13%  '
$ws.Range('D14').Value = '1.805.09'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').Value = '''0.643'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '34.389.96'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').Value = '''4.27'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '''69.35'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = '''255.21'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('D20').Value = '0.0₃0746'
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').Value = '''10.42'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').Value = '''4.24'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.08%  '
$ws.Range('E24').Value = '  -4.16%  '
$ws.Range('D25').Value = '''157.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').Value = '''16.43'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.69%  '
$ws.Range('D27').Value = '''7.07'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('E28').Value = '  -3.39%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  -2.79%  '
$ws.Range('D31').Value = '''0.0514'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.88%  '
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('D33').Value = '''3.59'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').Value = '''1.89'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.33%  '
$ws.Range('D35').Value = '1.453.77'
$ws.Range('E35').Value = '  -5.04%  '
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.633'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.0189'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''2.86'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.49%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '''83.03'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').Value = '''2.32'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('D42').Value = '''0.894'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('E43').Value = '  -2.25%  '
$ws.Range('D44').Value = '''0.0508'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.53%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''5.89'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '''1.04'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('D47').Value = '1.948.40'
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').Value = '''12.10'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').Value = '''99.40'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = '''49.64'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.73%  '
